# Apply updated cryptocurrency price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.182.47"
$ws.Range("E2").Value = "  +0.14%  "

$ws.Range("D3").Value = "2.052.04"
$ws.Range("E3").Value = "  -1.16%  "

$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").Value = "'248.25"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.02%  "

$ws.Range("E6").Value = "  -1.90%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").Value = "'56.85"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -3.96%  "

$ws.Range("D9").Value = "'0.382"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -2.69%  "

$ws.Range("E10").Value = "  -2.50%  "

$ws.Range("E11").Value = "  -0.05%  "

$ws.Range("D12").Value = "'16.25"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.55%  "

$ws.Range("D13").Value = "'0.883"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +7.75%  "

$ws.Range("D14").Value = "2.351.65"
$ws.Range("E14").Value = "  -1.17%  "

$ws.Range("E15").Value = "  +3.51%  "

$ws.Range("D16").Value = "2.067.77"
$ws.Range("E16").Value = "  -0.71%  "

$ws.Range("D17").Value = "'18.42"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +16.26%  "

$ws.Range("D18").Value = "37.177.33"
$ws.Range("E18").Value = "  -0.62%  "

$ws.Range("D19").Value = "'74.74"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.16%  "

$ws.Range("D20").Value = "0.0₃0897"
$ws.Range("E20").Value = "  -3.70%  "

$ws.Range("E21").Value = "  -1.41%  "

$ws.Range("D22").Value = "'237.23"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.94%  "

$ws.Range("E23").Value = "  -0.02%  "

$ws.Range("D24").Value = "'2.48"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +3.01%  "

$ws.Range("E25").Value = "  +2.01%  "

$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "'169.74"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.40%  "

$ws.Range("B27").Value = "PancakeSwap"
$ws.Range("C27").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D27").Value = "'2.17"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -4.80%  "

$ws.Range("D28").Value = "'20.09"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.51%  "

$ws.Range("E29").Value = "  -1.31%  "

$ws.Range("E30").Value = "  -1.31%  "

$ws.Range("E31").Value = "  +1.53%  "

$ws.Range("E32").Value = "  -2.47%  "

$ws.Range("E33").Value = "  -0.43%  "

$ws.Range("D34").Value = "'0.0891"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.07%  "

$ws.Range("E35").Value = "  +0.01%  "

$ws.Range("E36").Value = "  -2.73%  "

$ws.Range("E37").Value = "  +0.27%  "

$ws.Range("E38").Value = "  -2.35%  "

$ws.Range("D39").Value = "'5.29"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +14.75%  "

$ws.Range("E40").Value = "  +9.76%  "

$ws.Range("D41").Value = "'0.101"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -13.86%  "

$ws.Range("E42").Value = "  -2.06%  "

$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").Value = "'1.15"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.99%  "

$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").Value = "'17.30"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -3.34%  "

$ws.Range("D45").Value = "'96.09"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -3.22%  "

$ws.Range("D46").Value = "'2.43"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.81%  "

$ws.Range("D47").Value = "1.269.32"
$ws.Range("E47").Value = "  -2.91%  "

$ws.Range("E48").Value = "  -3.26%  "

$ws.Range("D49").Value = "'6.82"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.79%  "

$ws.Range("D50").Value = "2.241.25"
$ws.Range("E50").Value = "  -0.99%  "

$ws.Range("E51").Value = "  -0.94%  "
